$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")
$ws.Activate()

# 1) Fix the SQL query text in cell B3 (DiagnosisTab row): capitalize
#    "diagnosis" -> "Diagnosis" in the "Age at diagnosis (days)" column alias.
$oldText = $ws.Range("B3").Value2
$newText = $oldText -replace "Age at diagnosis \(days\)", "Age at Diagnosis (days)"
$ws.Range("B3").Value2 = $newText

# 2) Update the sheet view: scroll/select column C so the selection becomes
#    the whole column C (C1:C1048576) with the active cell at the top of the
#    visible area.
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
[void]$ws.Columns("C:C").Select()
